# Replace the three rows of shared-string labels (a/b, aa/bb, aaa/bbb) with
# 16 rows of plain numeric data in columns A and B, where each row's value
# equals its row number (A<n> = n, B<n> = n).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 1; $i -le 16; $i++) {
    $ws.Cells.Item($i, 1).Value = $i
    $ws.Cells.Item($i, 2).Value = $i
}

# Move the active selection to D9 (was E8).
[void]$ws.Range("D9").Select()
